# This workbook's sheet originally had:
#   Row 1  -> bold header labels (Lg., mm / Threading / .../ thread_size / material_surface)
#   Row 2+ -> data rows (group headers + fastener specs)
#
# The target state inserts a brand-new row above everything, filled with the
# plain numeric column indices 0-11 (styled like the old header row), pushes
# all existing rows down by one, and the row that used to be the header row
# (now row 2) loses its bold/bordered formatting as well as the values that
# used to live in the last two columns (K/L - "thread_size" / "material_surface").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shift every existing row down by one to make room for the new row 1.
$ws.Rows.Item(1).Insert()

# 2. The old header row (now row 2) keeps the bold/border/center-top style and
#    the K/L text from before the insert - strip that back down to the plain
#    "Normal" style and blank out K2/L2 so it matches the rest of the data rows.
$ws.Range("A2:L2").Style = "Normal"
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()

# 3. Fill the brand-new row 1 with the numeric column indices (0-11) and give
#    it the same look the header row used to have: bold, thin box border,
#    centered horizontally and top-aligned vertically.
for ($col = 1; $col -le 12; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

$headerRange = $ws.Range("A1:L1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.Item(1).LineStyle = 1
$headerRange.Borders.Item(2).LineStyle = 1
$headerRange.Borders.Item(3).LineStyle = 1
$headerRange.Borders.Item(4).LineStyle = 1
